$d = $word.ActiveDocument

# Update header date.
$d.Paragraphs.Item(1).Range.Text = "2026-02-25 Wednesday"

# Row 1 of the division-practice table.
$d.Paragraphs.Item(2).Range.Text = "89÷3="
$d.Paragraphs.Item(3).Range.Text = "24÷8="
$d.Paragraphs.Item(4).Range.Text = "88÷8="
$d.Paragraphs.Item(5).Range.Text = "29÷5="
$d.Paragraphs.Item(6).Range.Text = "23÷9="

# Row 2 (note: original 75÷8= cell is left in place; only its position
# among the five cells and the other four values change).
$d.Paragraphs.Item(26).Range.Text = "38÷3="
$d.Paragraphs.Item(27).Range.Text = "75÷8="
$d.Paragraphs.Item(28).Range.Text = "58÷6="
$d.Paragraphs.Item(29).Range.Text = "84÷8="
$d.Paragraphs.Item(30).Range.Text = "36÷7="

# Row 3.
$d.Paragraphs.Item(50).Range.Text = "49÷3="
$d.Paragraphs.Item(51).Range.Text = "95÷6="
$d.Paragraphs.Item(52).Range.Text = "23÷5="
$d.Paragraphs.Item(53).Range.Text = "25÷9="
$d.Paragraphs.Item(54).Range.Text = "19÷3="

# Row 4.
$d.Paragraphs.Item(74).Range.Text = "56÷4="
$d.Paragraphs.Item(75).Range.Text = "41÷6="
$d.Paragraphs.Item(76).Range.Text = "30÷6="
$d.Paragraphs.Item(77).Range.Text = "42÷6="
$d.Paragraphs.Item(78).Range.Text = "15÷8="

# Row 5.
$d.Paragraphs.Item(98).Range.Text = "25÷6="
$d.Paragraphs.Item(99).Range.Text = "86÷8="
$d.Paragraphs.Item(100).Range.Text = "91÷4="
$d.Paragraphs.Item(101).Range.Text = "79÷8="
$d.Paragraphs.Item(102).Range.Text = "49÷8="
